# Append new listings and refresh timestamps for the Lancers sheet.
# Final state: 13 data rows (2..14), each cell written explicitly so
# hyperlink relationship ids stay correctly paired with their F-column
# cell (row-insert style shifting does not renumber hyperlink refs here).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$rows = @(
    @('2025-11-15 12:32:28', '急募 【Python/AI/GAS 開発者・PM向け】「業務委託・再委託」の経験に関する30分インタビュー', 'システム開発', '1,000 ~ 5,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434693', 530, '🔥AI,Python ◆開発'),
    @('2025-11-15 12:32:28', '生成AI使用可 【急募】AIを活用したメールの自動応答システム構築依頼', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434695', 318, '🔥AI,Ai'),
    @('2025-11-15 12:32:28', '【急募】生成AI・RAG活用の業務ナレッジ検索システム改善', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434552', 318, '🔥AI,Ai'),
    @('2025-11-15 12:32:28', '海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)', 'システム開発', '5,000 円 ~ 10,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5251319', 135, '◆ツール,スクレイピング ◇サイト'),
    @('2025-11-15 12:32:28', '【システム開発】FileMaker Proを活用した販売システム構築', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434428', 118, '◆開発,システム開発'),
    @('2025-11-15 12:32:28', '初回 あるサイトの自動操作スクリプト開発(作業見積5時間以内/予算1万以内)の仕事・依頼', 'システム開発', '10,000 円 ~ 20,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434568', 85, '◆開発 ◇サイト'),
    @('2025-11-15 12:32:28', 'Flutter iOSアプリにおけるRevenueCat導入のバグ修正依頼', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434437', 38, '◇アプリ'),
    @('2025-11-15 12:32:28', '【自動応答】メール勧誘対策システムの構築依頼', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434698', 33, $null),
    @('2025-11-15 12:32:28', '【急募】料理教室のレシピデジタル化とマイページ構築依頼', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434648', 18, $null),
    @('2025-11-15 12:32:28', '【急募】TradingViewインジシグナルを用いたXAUUSD自動売買EA制作', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434524', 18, $null),
    @('2025-11-15 12:32:28', 'URL付きPDF資料の閲覧状況を可視化し、トラッキングする', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434431', 18, $null),
    @('2025-11-15 12:32:28', 'GAS構築できる方求む!', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434226', 13, $null),
    @('2025-11-15 12:32:28', '【Stable Diffusion】参考動画に沿って約100プロンプト構築', 'システム開発', '5,000 円 ~ 10,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5432055', 10, $null)
)

# Clear any previously-existing hyperlinks before we rewrite the grid,
# so stale ref/relationship pairs cannot linger.
$ws.Hyperlinks.Delete()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[5])
    $ws.Cells.Item($r, 7).Value = $row[6]
    if ($row[7] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[7]
    }
    $r = $r + 1
}

# Column B widened from 51 to 55 characters (ColumnWidth reads/writes
# with a fixed ~0.83-character offset vs. the stored OOXML width).
$ws.Columns.Item(2).ColumnWidth = 54.17
